$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column C (rows 2-16)
$cValues = @{
    2  = 216
    3  = 96
    4  = 89
    5  = 39
    6  = 215
    7  = 61
    8  = 54
    9  = 1528
    10 = 51
    11 = 92
    12 = 72
    13 = 38
    14 = 794
    15 = 40
    16 = 41
}

# New values for column D (rows 2-16)
$dValues = @{
    2  = 203.5
    3  = 85.5
    4  = 74.5
    5  = 21.5
    6  = 203.5
    7  = 47
    8  = 37.5
    9  = 1536.5
    10 = 31.5
    11 = 76.5
    12 = 58.5
    13 = 38
    14 = 855.5
    15 = 22.5
    16 = 23.5
}

foreach ($row in $cValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $cValues[$row]
}

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

# Update the summary average in C17
$ws.Cells.Item(17, 3).Value = 228.4
